$d = $word.ActiveDocument

# The document contains three "<id>...</id>" tag sequences, each split
# across three runs: "<id>" (Courier New, colored), the bare id text
# (plain formatting), and "</id>" (Courier New, colored). Each triple
# must be collapsed into a single run - "<id>p142v_N</id>" - carrying
# the formatting of the original "<id>" run.
#
# Approach: find the opening "<id>" run's text (left untouched, so its
# run properties/rsids survive unmodified), delete the text belonging to
# the two runs that follow it ("p142v_N" + "</id>"), then insert that
# same literal text right after the "<id>" run so Word appends it using
# that run's formatting, merging everything into one run.
#
# A single Range object is reused (and re-searched/collapsed) across the
# loop so each Find.Execute call resumes searching after the previous
# match instead of restarting from the top of the document.

$ids = @("p142v_1", "p142v_2", "p142v_3")

$rng = $d.Content
foreach ($id in $ids) {
    $found = $rng.Find.Execute("<id>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

    $tail = $id + "</id>"
    $rng2 = $d.Range($rng.End, $rng.End + $tail.Length)
    $rng2.Delete()

    $rng.InsertAfter($tail)

    # Move past this match so the next Find.Execute looks further ahead.
    $rng.Collapse(0)
}
